$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = -22.15180000000002
$ws.Range("D18").Value = -8.644199999999998
$ws.Range("A21").Value = -20.15059999999998
$ws.Range("A23").Value = -20.14029999999998
$ws.Range("B24").Value = 5.600100000000003
$ws.Range("A25").Value = -21.79919999999999
$ws.Range("B28").Value = 6.045700000000004
$ws.Range("B36").Value = 9.498100000000004
$ws.Range("B45").Value = 5.128500000000005
$ws.Range("B48").Value = 6.868400000000006
$ws.Range("B49").Value = 5.957199999999998
$ws.Range("D51").Value = -7.928599999999999
$ws.Range("B52").Value = 5.5224
$ws.Range("A53").Value = -21.76239999999999
$ws.Range("B53").Value = 5.836899999999999
$ws.Range("B54").Value = 4.796700000000005
$ws.Range("D55").Value = -8.942400000000001
$ws.Range("A57").Value = -22.30770000000001
$ws.Range("A59").Value = -22.2779
$ws.Range("D64").Value = -7.400299999999991
$ws.Range("A69").Value = -21.65199999999999
$ws.Range("B70").Value = 7.115300000000004
$ws.Range("A79").Value = -20.09280000000001
$ws.Range("D80").Value = -8.044899999999998
$ws.Range("A83").Value = -21.76839999999999
$ws.Range("B86").Value = 5.057800000000001
$ws.Range("B87").Value = 5.340600000000002
$ws.Range("D92").Value = -6.8481
$ws.Range("A93").Value = -21.48880000000002
$ws.Range("D94").Value = -6.394099999999998
$ws.Range("D96").Value = -8.370900000000002
$ws.Range("B101").Value = 4.6294
